# Update IFRS financial figures for 빙그레 (company_list sheet, rows 2-9)
# per corrected source data ("error solve ifrs list").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 8200
$ws.Range("E2").Value = 417
$ws.Range("F2").Value = 418
$ws.Range("G2").Value = 485
$ws.Range("H2").Value = 376
$ws.Range("I2").Value = 379
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 5789
$ws.Range("L2").Value = 948
$ws.Range("M2").Value = 4840
$ws.Range("N2").Value = 4840
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 498
$ws.Range("Q2").Value = 444
$ws.Range("R2").Value = -533
$ws.Range("S2").Value = -122
$ws.Range("T2").Value = 370
$ws.Range("U2").Value = 74
$ws.Range("V2").Value = 28
$ws.Range("W2").Value = 5.08
$ws.Range("X2").Value = 4.58
$ws.Range("Y2").Value = 7.98
$ws.Range("Z2").Value = 6.57
$ws.Range("AA2").Value = 19.6
$ws.Range("AB2").Value = 902.63
$ws.Range("AC2").Value = 3843
$ws.Range("AD2").Value = 19.62
$ws.Range("AE2").Value = 54741
$ws.Range("AF2").Value = 1.38
$ws.Range("AG2").Value = 1250
$ws.Range("AH2").Value = 1.66
$ws.Range("AI2").Value = 29.19
$ws.Range("AJ2").Value = 9851241

# Row 3
$ws.Range("D3").Value = 7996
$ws.Range("E3").Value = 317
$ws.Range("F3").Value = 317
$ws.Range("G3").Value = 324
$ws.Range("H3").Value = 247
$ws.Range("I3").Value = 247
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6170
$ws.Range("L3").Value = 998
$ws.Range("M3").Value = 5172
$ws.Range("N3").Value = 5172
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 498
$ws.Range("Q3").Value = 647
$ws.Range("R3").Value = -419
$ws.Range("S3").Value = -138
$ws.Range("T3").Value = 260
$ws.Range("U3").Value = 387
$ws.Range("V3").Value = 1
$ws.Range("W3").Value = 3.96
$ws.Range("X3").Value = 3.09
$ws.Range("Y3").Value = 4.94
$ws.Range("Z3").Value = 4.14
$ws.Range("AA3").Value = 19.3
$ws.Range("AB3").Value = 928.98
$ws.Range("AC3").Value = 2512
$ws.Range("AD3").Value = 27.27
$ws.Range("AE3").Value = 58490
$ws.Range("AF3").Value = 1.17
$ws.Range("AG3").Value = 1250
$ws.Range("AH3").Value = 1.82
$ws.Range("AI3").Value = 44.67
$ws.Range("AJ3").Value = 9851241

# Row 4
$ws.Range("D4").Value = 8132
$ws.Range("E4").Value = 372
$ws.Range("F4").Value = 372
$ws.Range("G4").Value = 370
$ws.Range("H4").Value = 287
$ws.Range("I4").Value = 287
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6238
$ws.Range("L4").Value = 1020
$ws.Range("M4").Value = 5218
$ws.Range("N4").Value = 5218
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 498
$ws.Range("Q4").Value = 656
$ws.Range("R4").Value = -500
$ws.Range("S4").Value = -112
$ws.Range("T4").Value = 260
$ws.Range("U4").Value = 396
$ws.Range("V4").ClearContents()
$ws.Range("W4").Value = 4.58
$ws.Range("X4").Value = 3.53
$ws.Range("Y4").Value = 5.53
$ws.Range("Z4").Value = 4.63
$ws.Range("AA4").Value = 19.55
$ws.Range("AB4").Value = 967.75
$ws.Range("AC4").Value = 2916
$ws.Range("AD4").Value = 21.88
$ws.Range("AE4").Value = 59010
$ws.Range("AF4").Value = 1.08
$ws.Range("AG4").Value = 1250
$ws.Range("AH4").Value = 1.96
$ws.Range("AI4").Value = 38.47
$ws.Range("AJ4").Value = 9851241

# Row 5
$ws.Range("D5").Value = 8147
$ws.Range("E5").Value = 347
$ws.Range("F5").Value = 347
$ws.Range("G5").Value = 481
$ws.Range("H5").Value = 296
$ws.Range("I5").Value = 296
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 6356
$ws.Range("L5").Value = 1068
$ws.Range("M5").Value = 5288
$ws.Range("N5").Value = 5288
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 498
$ws.Range("Q5").Value = 419
$ws.Range("R5").Value = -284
$ws.Range("S5").Value = -111
$ws.Range("T5").Value = 230
$ws.Range("U5").Value = 188
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 4.26
$ws.Range("X5").Value = 3.64
$ws.Range("Y5").Value = 5.64
$ws.Range("Z5").Value = 4.7
$ws.Range("AA5").Value = 20.19
$ws.Range("AB5").Value = 1003.32
$ws.Range("AC5").Value = 3006
$ws.Range("AD5").Value = 20.99
$ws.Range("AE5").Value = 59812
$ws.Range("AF5").Value = 1.05
$ws.Range("AG5").Value = 1250
$ws.Range("AH5").Value = 1.98
$ws.Range("AI5").Value = 37.32
$ws.Range("AJ5").Value = 9851241

# Row 6
$ws.Range("D6").Value = 8552
$ws.Range("E6").Value = 393
$ws.Range("F6").Value = 393
$ws.Range("G6").Value = 444
$ws.Range("H6").Value = 342
$ws.Range("I6").Value = 342
$ws.Range("K6").Value = 6539
$ws.Range("L6").Value = 1089
$ws.Range("M6").Value = 5451
$ws.Range("N6").Value = 5451
$ws.Range("P6").Value = 498
$ws.Range("Q6").Value = 726
$ws.Range("R6").Value = -515
$ws.Range("S6").Value = -117
$ws.Range("T6").Value = 211
$ws.Range("U6").Value = 515
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 4.6
$ws.Range("X6").Value = 4
$ws.Range("Y6").Value = 6.37
$ws.Range("Z6").Value = 5.3
$ws.Range("AA6").Value = 19.97
$ws.Range("AB6").Value = 1043.83
$ws.Range("AC6").Value = 3472
$ws.Range("AD6").Value = 20.88
$ws.Range("AE6").Value = 61647
$ws.Range("AF6").Value = 1.18
$ws.Range("AG6").Value = 1350
$ws.Range("AH6").Value = 1.86
$ws.Range("AI6").Value = 34.9
$ws.Range("AJ6").Value = 9851241

# Row 7
$ws.Range("D7").Value = 8737
$ws.Range("E7").Value = 435
$ws.Range("G7").Value = 493
$ws.Range("H7").Value = 380
$ws.Range("I7").Value = 380
$ws.Range("K7").Value = 6817
$ws.Range("L7").Value = 1107
$ws.Range("M7").Value = 5710
$ws.Range("N7").Value = 5710
$ws.Range("P7").Value = 500
$ws.Range("Q7").Value = 693
$ws.Range("R7").Value = -535
$ws.Range("S7").Value = -247
$ws.Range("T7").Value = 343
$ws.Range("U7").Value = 550
$ws.Range("W7").Value = 4.98
$ws.Range("X7").Value = 4.35
$ws.Range("Y7").Value = 6.81
$ws.Range("Z7").Value = 5.69
$ws.Range("AA7").Value = 19.38
$ws.Range("AC7").Value = 3857
$ws.Range("AD7").Value = 13.71
$ws.Range("AE7").Value = 64580
$ws.Range("AF7").Value = 0.82
$ws.Range("AG7").Value = 1400
$ws.Range("AH7").Value = 2.65
$ws.Range("AI7").Value = 36.29

# Row 8
$ws.Range("D8").Value = 9069
$ws.Range("E8").Value = 499
$ws.Range("G8").Value = 550
$ws.Range("H8").Value = 420
$ws.Range("I8").Value = 420
$ws.Range("K8").Value = 7137
$ws.Range("L8").Value = 1130
$ws.Range("M8").Value = 6003
$ws.Range("N8").Value = 6003
$ws.Range("P8").Value = 500
$ws.Range("Q8").Value = 607
$ws.Range("R8").Value = -457
$ws.Range("S8").Value = -103
$ws.Range("T8").Value = 310
$ws.Range("U8").Value = 310
$ws.Range("W8").Value = 5.5
$ws.Range("X8").Value = 4.63
$ws.Range("Y8").Value = 7.17
$ws.Range("Z8").Value = 6.02
$ws.Range("AA8").Value = 18.82
$ws.Range("AC8").Value = 4263
$ws.Range("AD8").Value = 12.41
$ws.Range("AE8").Value = 67897
$ws.Range("AF8").Value = 0.78
$ws.Range("AG8").Value = 1462
$ws.Range("AH8").Value = 2.76
$ws.Range("AI8").Value = 34.3

# Row 9
$ws.Range("D9").Value = 9353
$ws.Range("E9").Value = 537
$ws.Range("G9").Value = 597
$ws.Range("H9").Value = 460
$ws.Range("I9").Value = 460
$ws.Range("K9").Value = 7460
$ws.Range("L9").Value = 1133
$ws.Range("M9").Value = 6323
$ws.Range("N9").Value = 6323
$ws.Range("P9").Value = 500
$ws.Range("Q9").Value = 613
$ws.Range("R9").Value = -453
$ws.Range("S9").Value = -127
$ws.Range("T9").Value = 297
$ws.Range("U9").Value = 320
$ws.Range("W9").Value = 5.74
$ws.Range("X9").Value = 4.92
$ws.Range("Y9").Value = 7.46
$ws.Range("Z9").Value = 6.3
$ws.Range("AA9").Value = 17.92
$ws.Range("AC9").Value = 4669
$ws.Range("AD9").Value = 11.33
$ws.Range("AE9").Value = 71516
$ws.Range("AF9").Value = 0.74
$ws.Range("AG9").Value = 1488
$ws.Range("AH9").Value = 2.81
$ws.Range("AI9").Value = 31.86
